# Backup: Save all current work including dashboard generations and README updates
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) STRATEGY_SELECTOR: add "SECTION D: STRATEGY ALIGNMENT" block (rows 25-29)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("STRATEGY_SELECTOR")

# Section header - reuse the existing "SECTION x:" style from A18
$ws2.Range("A18").Copy()
$ws2.Range("A25").PasteSpecial(-4122)
$ws2.Range("A25").Value = "SECTION D: STRATEGY ALIGNMENT"

# Table header row (Target / Goal / Current / Status)
$ws2.Range("A26").Value = "Target"
$ws2.Range("B26").Value = "Goal"
$ws2.Range("C26").Value = "Current"
$ws2.Range("D26").Value = "Status"
$hdr = $ws2.Range("A26:D26")
$hdr.Font.Bold = $true
$hdr.Font.Color = 16777215
$hdr.Interior.Color = 7954721
$hdr.Borders.LineStyle = 1

# Row 27: Carbon Neutral Year
$ws2.Range("A27").Value = "Carbon Neutral Year"
$ws2.Range("B27").Value = 2030
$ws2.Range("C27").Value = "TBD"
$ws2.Range("D27").Formula = '=IF(C27="TBD", "PENDING", IF(C27<=B27, "ON TRACK", "BEHIND"))'

# Row 28: Product Recyclability %
$ws2.Range("A28").Value = "Product Recyclability %"
$ws2.Range("B28").Value = 1
$ws2.Range("C28").Value = 0
$ws2.Range("D28").Formula = '=IF(C28>=B28, "COMPLIANT", "GAP: "&TEXT(B28-C28,"0%"))'

# Row 29: Fair Wage (vs Market)
$ws2.Range("A29").Value = "Fair Wage (vs Market)"
$ws2.Range("B29").Value = 1.1
$ws2.Range("C29").Value = 1
$ws2.Range("D29").Formula = '=IF(C29>=B29, "COMPLIANT", "RAISE WAGES")'

# Style A27:A29 + D27:D29 - plain bordered cells (reuse existing border style)
$ws2.Range("A20").Copy()
$left = $ws2.Range("A27:A29")
$left.PasteSpecialFormatsOnly = $true
$ws2.Range("A20").Copy()
$ws2.Range("A27:A29").PasteSpecial(-4122)
$ws2.Range("A20").Copy()
$ws2.Range("D27:D29").PasteSpecial(-4122)
$ws2.Range("A27:A29").Borders.LineStyle = 1
$ws2.Range("D27:D29").Borders.LineStyle = 1

# C27: reuse the yellow "quantity" input style (same as B13 etc.)
$ws2.Range("B13").Copy()
$ws2.Range("C27").PasteSpecial(-4122)
$ws2.Range("C27").Value = "TBD"

# B28: percent format with border (no fill) - new style
$ws2.Range("B28").NumberFormat = "0%"
$ws2.Range("B28").Borders.LineStyle = 1

# C28: reuse existing percent-with-fill style (same as B15 row - s=21)
$ws2.Range("B15").Copy()
$ws2.Range("C28").PasteSpecial(-4122)
$ws2.Range("C28").Value = 0

# B29 / C29: new "0.0x" multiple format
$ws2.Range("B29").NumberFormat = "0.0x"
$ws2.Range("B29").Borders.LineStyle = 1
$ws2.Range("C29").NumberFormat = "0.0x"
$ws2.Range("C29").Interior.Color = 13434879
$ws2.Range("C29").Borders.LineStyle = 1

# Conditional formatting on D27:D29 - green for ON TRACK / COMPLIANT,
# red/pink for GAP / BEHIND / RAISE WAGES
$rng = $ws2.Range("D27:D29")
$fcGood = $rng.FormatConditions.Add(2, 0, 'OR(D27="ON TRACK", D27="COMPLIANT")')
$fcGood.Interior.Color = 13561798
$fcBad = $rng.FormatConditions.Add(2, 0, 'OR(LEFT(D27,3)="GAP", D27="BEHIND", D27="RAISE WAGES")')
$fcBad.Interior.Color = 13551615

# ---------------------------------------------------------------------------
# 2) Add the new UPLOAD_READY_ESG sheet at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "UPLOAD_READY_ESG"
$ws3.Outline.SummaryRow = 1
$ws3.Outline.SummaryColumn = 1

$ws3.Range("A1").Value = "UPLOAD READY DATA - DO NOT EDIT"
$ws3.Range("A1").Font.Bold = $true
$ws3.Range("A1").Font.Color = 255

$ws3.Range("A5").Value = "Initiative"
$ws3.Range("B5").Value = "Quantity"
$ws3.Range("C5").Value = "Investment"
$ws3.Range("D5").Value = "CO2 Reduction"
$ws3.Range("E5").Value = "Tax Savings"
$hdr3 = $ws3.Range("A5:E5")
$hdr3.Font.Bold = $true
$hdr3.Borders.LineStyle = 1

$ws3.Range("A6").Value = "Solar PV Panels"
$ws3.Range("B6").Formula = "=STRATEGY_SELECTOR!B13"
$ws3.Range("C6").Formula = "=STRATEGY_SELECTOR!C13"
$ws3.Range("D6").Formula = "=STRATEGY_SELECTOR!D13"
$ws3.Range("E6").Formula = "=STRATEGY_SELECTOR!E13"

$ws3.Range("A7").Value = "Trees Planted"
$ws3.Range("B7").Formula = "=STRATEGY_SELECTOR!B14"
$ws3.Range("C7").Formula = "=STRATEGY_SELECTOR!C14"
$ws3.Range("D7").Formula = "=STRATEGY_SELECTOR!D14"
$ws3.Range("E7").Formula = "=STRATEGY_SELECTOR!E14"

$ws3.Range("A8").Value = "Green Electricity"
$ws3.Range("B8").Formula = "=STRATEGY_SELECTOR!B15"
$ws3.Range("C8").Formula = "=STRATEGY_SELECTOR!C15"
$ws3.Range("D8").Formula = "=STRATEGY_SELECTOR!D15"
$ws3.Range("E8").Formula = "=STRATEGY_SELECTOR!E15"

$ws3.Range("A9").Value = "CO2 Credits"
$ws3.Range("B9").Formula = "=STRATEGY_SELECTOR!B16"
$ws3.Range("C9").Formula = "=STRATEGY_SELECTOR!C16"
$ws3.Range("D9").Formula = "=STRATEGY_SELECTOR!D16"
$ws3.Range("E9").Formula = "=STRATEGY_SELECTOR!E16"

# Strip the auto-inherited number formats picked up from STRATEGY_SELECTOR
$ws3.Range("B6:E9").Style = "Normal"

$ws3.Columns.Item(1).ColumnWidth = 19.166666666666668
$ws3.Columns.Item(2).ColumnWidth = 14.166666666666666
$ws3.Columns.Item(3).ColumnWidth = 14.166666666666666
$ws3.Columns.Item(4).ColumnWidth = 14.166666666666666
$ws3.Columns.Item(5).ColumnWidth = 14.166666666666666

Write-Output "Applied ESG dashboard backup changes"
